$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("G3").Value = 80
$ws.Range("F5").Value = 5463
$ws.Range("G5").Value = 80
$ws.Range("F6").Value = 5463
$ws.Range("G6").Value = 80
$ws.Range("F7").Value = 197
$ws.Range("F14").Value = 794
$ws.Range("F15").Value = 6456
$ws.Range("F16").Value = 37
$ws.Range("F17").Value = 84
$ws.Range("F19").Value = 4201
$ws.Range("F22").Value = 4138
$ws.Range("F23").Value = 4058
$ws.Range("F24").Value = 197
$ws.Range("F25").Value = 199
$ws.Range("F26").Value = 270
$ws.Range("F30").Value = 148
$ws.Range("F31").Value = 38
$ws.Range("F33").Value = 143
$ws.Range("F34").Value = 51
$ws.Range("F35").Value = 7292
$ws.Range("F37").Value = 1212
$ws.Range("F38").Value = 603
$ws.Range("F40").Value = 978
$ws.Range("F42").Value = 1460
$ws.Range("F43").Value = 190
$ws.Range("F44").Value = 801
$ws.Range("F45").Value = 33
$ws.Range("F46").Value = 3498
$ws.Range("F47").Value = 330
$ws.Range("F49").Value = 807
$ws.Range("F50").Value = 1010

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 17
$ws.Range("F5").Value = 98

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G3").Value = 80
$ws.Range("F5").Value = 5463
$ws.Range("G5").Value = 80
$ws.Range("F6").Value = 5463
$ws.Range("G6").Value = 80
$ws.Range("F7").Value = 197
$ws.Range("F9").Value = 98
$ws.Range("F15").Value = 794
$ws.Range("F16").Value = 6456
$ws.Range("F17").Value = 37
$ws.Range("F18").Value = 84
$ws.Range("F20").Value = 4201
$ws.Range("F23").Value = 4138
$ws.Range("F24").Value = 4058
$ws.Range("F25").Value = 197
$ws.Range("F26").Value = 199
$ws.Range("F27").Value = 270
$ws.Range("F33").Value = 7292
$ws.Range("F35").Value = 1212
$ws.Range("F36").Value = 603
$ws.Range("F39").Value = 978
$ws.Range("F41").Value = 1460
$ws.Range("F42").Value = 190
$ws.Range("F43").Value = 801
$ws.Range("F44").Value = 33
$ws.Range("F45").Value = 3498
$ws.Range("F46").Value = 330
$ws.Range("F48").Value = 807
$ws.Range("F49").Value = 1010
